$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 = "Enterprises density (per 1000 people)"
# C11 (SMEs column) : 1.8  -> 1.78
# D11 (MSMEs column) : 11.7 -> 11.68
# Force text format so the numeric-looking values are kept as text,
# matching the shared-string text cells in the workbook.
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "1.78"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "11.68"
